$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting N:P -> O:Q
$ws.Columns("N").Insert()

# Update the selected cell/range on the sheet
$ws.Range("K17").Select()
